# Fruta / hortaliza, semanal
# Insert a new weekly record at row 216 (Vega Monumental Concepción - Piña),
# pushing the existing rows 216-307 down to 217-308.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 216 - shifts rows 216:307 -> 217:308
# and extends the used range to A1:T308.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(216, 1).Value  = 11
$ws.Cells.Item(216, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(216, 3).Value  = "Bíobío"
$ws.Cells.Item(216, 4).Value  = 45205
$ws.Cells.Item(216, 5).Value  = 8
$ws.Cells.Item(216, 6).Value  = "Fruta"
$ws.Cells.Item(216, 7).Value  = 100108
$ws.Cells.Item(216, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(216, 9).Value  = 100108005
$ws.Cells.Item(216, 10).Value = "Piña"
$ws.Cells.Item(216, 11).Value = "Caramelo"
$ws.Cells.Item(216, 12).Value = "Segunda"
$ws.Cells.Item(216, 13).Value = 230
$ws.Cells.Item(216, 14).Value = 21000
$ws.Cells.Item(216, 15).Value = 21000
$ws.Cells.Item(216, 16).Value = 21000
$ws.Cells.Item(216, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(216, 18).Value = "Ecuador"
$ws.Cells.Item(216, 19).Value = 1500
$ws.Cells.Item(216, 20).Value = 14
